# Regenerate s_val data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the derived G (sum)
# column for each data row (2-28) on the active worksheet, leaving
# column A (date) and column F (Win) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    3  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    4  = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.242251378316819)
    5  = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    6  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    7  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    8  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    9  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    10 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    11 = @(0.1190320826869504, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0.8034070775528621)
    12 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    13 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086)
    14 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    15 = @(0.6606524410359556, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 2.214453472130288)
    16 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    17 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    18 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    19 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    20 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    21 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    22 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    23 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    24 = @(3.286832544864788, 117.745847958593, 3.537761648806719, 10.19245300693656, 134.7628951592011)
    25 = @(0.01293466051926884, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0.9634143985795411)
    26 = @(0.1190320826869504, 0.306821227259698, 22.3905356188092, 0.4942365360607697, 23.31062546481661)
    27 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302)
    28 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
